$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.071837902069092
$ws.Range("B1").Value = 2.473607778549194
$ws.Range("C1").Value = 2.573086500167847
$ws.Range("D1").Value = 3.33726978302002
$ws.Range("E1").Value = 0.9749209880828857
